# -----------------------------------------------------------------------
# MAL-MMG_CUSTOM reference.xlsx update:
#   - bump Package info metadata (version/DHIS2 version/build/last updated/name)
#   - reorder + update dataElements rows (category combo "Sex" introduced)
#   - reorder dataElementGroups link rows to match
#   - insert a new "categoryCombos" sheet (between dataElementGroups and userGroups)
#   - update userGroups rows (new "Last updated" dates, reordered)
#   - users sheet is carried over unchanged
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Package info
# ---------------------------------------------------------------------
$pkg = $wb.Worksheets.Item("Package info")

$pkg.Cells.Item(4,2).Value = "1.2.0"
$pkg.Cells.Item(5,2).Value = "2.34.4"
$pkg.Cells.Item(6,1).Value = "DHIS2 build"
$pkg.Cells.Item(6,2).Value = "aff07fb"
$pkg.Cells.Item(7,1).Value = "Last updated"
$pkg.Cells.Item(7,2).Value = "20210520T090044"

# New row 8 - copy formatting from row 6 (same banding as row 8 would get) then set values
$pkg.Range("A6:B6").Copy($pkg.Range("A8:B8"))
$pkg.Cells.Item(8,1).Value = "Name"
$pkg.Cells.Item(8,2).Value = "MAL-MMG_CUSTOM_V1.2.0_2.34.4-en"

# Column B got narrower (49.71 -> 33.71 chars); engine quantizes ColumnWidth to
# 1/6 increments, so 32.8 is the closest achievable match to 33.7109375
$pkg.Columns.Item(2).ColumnWidth = 32.8

# ---------------------------------------------------------------------
# 2. dataElements - rows reordered + E3/E6 categorycombo renamed to "Sex"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("dataElements")

$de.Cells.Item(2,1).Value = "MAL - Migrant and mobile population (MMP) positive"
$de.Cells.Item(2,2).Value = "MMP positive"
$de.Cells.Item(2,3).Value = "MAL_MM_POP_POS"
$de.Cells.Item(2,4).Value = "Number of migrant and mobile population (MMP) that are positive with either microscopy and RDT"
$de.Cells.Item(2,5).Value = "bjDvmb4bfuf"
$de.Cells.Item(2,6).Value = "2019-10-20"
$de.Cells.Item(2,7).Value = "CWHBMa4nC9J"

$de.Cells.Item(3,1).Value = "MAL - Malaria tested from cross-borders"
$de.Cells.Item(3,2).Value = "Tested from cross-borders"
$de.Cells.Item(3,3).Value = "MAL_TEST_CROSS_BORDERS"
$de.Cells.Item(3,4).Value = "Number of suspected cases tested with either microscopy or RDT in administrative area bordering international border"
$de.Cells.Item(3,5).Value = "Sex"
$de.Cells.Item(3,6).Value = "2019-10-20"
$de.Cells.Item(3,7).Value = "CxI1FHE4oEh"

$de.Cells.Item(4,1).Value = "MAL - Migrant and mobile population (MMP) followed up for 14 days"
$de.Cells.Item(4,2).Value = "MMP followed up for 14 days"
$de.Cells.Item(4,3).Value = "MAL_MM_POP_FOLLO_UP_FOR_14D"
$de.Cells.Item(4,4).Value = "Number of migrant and mobile population (MMP) that are are followed-up for 14 days (with testing using microscopy or RDT at intervals)"
$de.Cells.Item(4,5).Value = "bjDvmb4bfuf"
$de.Cells.Item(4,6).Value = "2019-10-20"
$de.Cells.Item(4,7).Value = "kdMT3AuDzj1"

$de.Cells.Item(5,1).Value = "MAL - Migrant and mobile population (MMP) tested"
$de.Cells.Item(5,2).Value = "MMP tested"
$de.Cells.Item(5,3).Value = "MAL_MM_POP_TEST"
$de.Cells.Item(5,4).Value = "Number of migrant and mobile population (MMP) suspected and tested with either microscopy and RDT"
$de.Cells.Item(5,5).Value = "bjDvmb4bfuf"
$de.Cells.Item(5,6).Value = "2019-10-20"
$de.Cells.Item(5,7).Value = "S3AqkeU4DET"

$de.Cells.Item(6,1).Value = "MAL - Malaria positive from cross-borders"
$de.Cells.Item(6,2).Value = "Positive from cross-borders"
$de.Cells.Item(6,3).Value = "MAL_POS_CROSS_BORDERS"
$de.Cells.Item(6,4).Value = "Number of positive cases with either microscopy or RDT in administrative area bordering international border"
$de.Cells.Item(6,5).Value = "Sex"
$de.Cells.Item(6,6).Value = "2019-10-20"
$de.Cells.Item(6,7).Value = "UwaQ0MJzXBz"

$de.Cells.Item(7,1).Value = "MAL - Malaria positive from cross-borders followed for 14 days"
$de.Cells.Item(7,2).Value = "Positive from cross-borders followed for 14 days"
$de.Cells.Item(7,3).Value = "MAL_POS_CROSS_BORDERS_FOLLO_14D"
$de.Cells.Item(7,4).Value = "Number of suspected malaria cases positive with either microscopy or RDT in administrative area bordering international border"
$de.Cells.Item(7,5).Value = "bjDvmb4bfuf"
$de.Cells.Item(7,6).Value = "2019-10-20"
$de.Cells.Item(7,7).Value = "wAHUeGPbH9A"

# ---------------------------------------------------------------------
# 3. dataElementGroups - column B reordered to match dataElements reorder
# ---------------------------------------------------------------------
$deg = $wb.Worksheets.Item("dataElementGroups")

$deg.Cells.Item(2,2).Value = "MAL - Migrant and mobile population (MMP) positive"
$deg.Cells.Item(3,2).Value = "MAL - Malaria tested from cross-borders"
$deg.Cells.Item(4,2).Value = "MAL - Migrant and mobile population (MMP) followed up for 14 days"
$deg.Cells.Item(5,2).Value = "MAL - Migrant and mobile population (MMP) tested"
$deg.Cells.Item(6,2).Value = "MAL - Malaria positive from cross-borders"
$deg.Cells.Item(7,2).Value = "MAL - Malaria positive from cross-borders followed for 14 days"

# ---------------------------------------------------------------------
# 4. New sheet: categoryCombos (inserted after dataElementGroups)
# ---------------------------------------------------------------------
$after = $wb.Worksheets.Item("dataElementGroups")
$cc = $wb.Worksheets.Add($null, $after)
$cc.Name = "categoryCombos"

# Reuse existing banded-row styles (header / row) by copying formatted
# ranges from an existing sheet, then overwrite with the new text.
$de.Range("A1:D1").Copy($cc.Range("A1:D1"))
$de.Range("A2:D2").Copy($cc.Range("A2:D2"))

$cc.Cells.Item(1,1).Value = "Name"
$cc.Cells.Item(1,2).Value = "Last updated"
$cc.Cells.Item(1,3).Value = "UID"
$cc.Cells.Item(1,4).Value = "Categories"

$cc.Cells.Item(2,1).Value = "Sex"
$cc.Cells.Item(2,2).Value = "'2017-06-02"
$cc.Cells.Item(2,3).Value = "VkQPxB6VdoG"
$cc.Cells.Item(2,4).Value = " "

# Column widths (engine quantizes ColumnWidth to 1/6 increments, so these
# are the closest achievable values to 6.71 / 14.71 / 13.71 / 12.71 chars)
$cc.Columns.Item(1).ColumnWidth = 5.8
$cc.Columns.Item(2).ColumnWidth = 13.8
$cc.Columns.Item(3).ColumnWidth = 12.8
$cc.Columns.Item(4).ColumnWidth = 11.8

# ---------------------------------------------------------------------
# 5. userGroups - rows reordered + "Last updated" bumped to 2021-05-20
# ---------------------------------------------------------------------
$ug = $wb.Worksheets.Item("userGroups")

$ug.Cells.Item(2,1).Value = "Malaria data capture"
$ug.Cells.Item(2,2).Value = "'2021-05-20"
$ug.Cells.Item(2,3).Value = "fRSrUJ6SMGH"

$ug.Cells.Item(3,1).Value = "Malaria admin"
$ug.Cells.Item(3,2).Value = "'2021-05-20"
$ug.Cells.Item(3,3).Value = "suMb19wGXPR"

$ug.Cells.Item(4,1).Value = "Malaria access"
$ug.Cells.Item(4,2).Value = "'2021-05-20"
$ug.Cells.Item(4,3).Value = "ZXEVDM9XRea"

# ---------------------------------------------------------------------
# 6. users sheet is unchanged - leave as-is.
# ---------------------------------------------------------------------

# Restore the originally-active sheet selection.
$pkg.Activate()
